$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.276052666666667
$ws.Range("H2").Value = 6.828158
$ws.Range("I2").Value = 0.005247614157263819
$ws.Range("J2").Value = 0.005247614157263819
$ws.Range("M2").Value = 145.7007446666667
$ws.Range("N2").Value = 437.1022340000001
$ws.Range("O2").Value = 0.2865937750105843
$ws.Range("P2").Value = 0.2865937750105843
$ws.Range("Q2").Value = 331.6225684338858
$ws.Range("R2").Value = 2984.603115904973
$ws.Range("S2").Value = 0.001503933551129223
$ws.Range("T2").Value = 0.001503933551129224

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.276052666666667
$ws.Range("H3").Value = 6.828158
$ws.Range("I3").Value = 0.005247614157263819
$ws.Range("J3").Value = 0.005247614157263819
$ws.Range("O3").Value = 0.3320294904365841
$ws.Range("P3").Value = 0.3320294904365841
$ws.Range("Q3").Value = 384.1970133870072
$ws.Range("R3").Value = 3457.773120483065
$ws.Range("S3").Value = 0.001742362654644111
$ws.Range("T3").Value = 0.001742362654644111

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.276052666666667
$ws.Range("H4").Value = 6.828158
$ws.Range("I4").Value = 0.005247614157263819
$ws.Range("J4").Value = 0.005247614157263819
$ws.Range("M4").Value = 128.1261546666667
$ws.Range("N4").Value = 384.378464
$ws.Range("O4").Value = 0.2520245069956105
$ws.Range("P4").Value = 0.2520245069956105
$ws.Range("Q4").Value = 291.6218759988125
$ws.Range("R4").Value = 2624.596883989312
$ws.Range("S4").Value = 0.0013225273708876
$ws.Range("T4").Value = 0.0013225273708876

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.276052666666667
$ws.Range("H5").Value = 6.828158
$ws.Range("I5").Value = 0.005247614157263819
$ws.Range("J5").Value = 0.005247614157263819
$ws.Range("M5").Value = 65.761079
$ws.Range("N5").Value = 197.283237
$ws.Range("O5").Value = 0.1293522275572212
$ws.Range("P5").Value = 0.1293522275572212
$ws.Range("Q5").Value = 149.6756792208273
$ws.Range("R5").Value = 1347.081112987446
$ws.Range("S5").Value = 0.0006787905806028852
$ws.Range("T5").Value = 0.0006787905806028853

$ws.Range("I6").Value = 0.1062533062835484
$ws.Range("J6").Value = 0.1062533062835484
$ws.Range("M6").Value = 145.7007446666667
$ws.Range("N6").Value = 437.1022340000001
$ws.Range("O6").Value = 0.2865937750105843
$ws.Range("P6").Value = 0.2865937750105843
$ws.Range("Q6").Value = 6714.669424688649
$ws.Range("R6").Value = 60432.02482219784
$ws.Range("S6").Value = 0.03045153615515796
$ws.Range("T6").Value = 0.03045153615515797

$ws.Range("I7").Value = 0.1062533062835484
$ws.Range("J7").Value = 0.1062533062835484
$ws.Range("O7").Value = 0.3320294904365841
$ws.Range("P7").Value = 0.3320294904365841
$ws.Range("S7").Value = 0.03527923114252887
$ws.Range("T7").Value = 0.03527923114252888

$ws.Range("I8").Value = 0.1062533062835484
$ws.Range("J8").Value = 0.1062533062835484
$ws.Range("M8").Value = 128.1261546666667
$ws.Range("N8").Value = 384.378464
$ws.Range("O8").Value = 0.2520245069956105
$ws.Range("P8").Value = 0.2520245069956105
$ws.Range("Q8").Value = 5904.738340297722
$ws.Range("R8").Value = 53142.64506267949
$ws.Range("S8").Value = 0.02677843713276488
$ws.Range("T8").Value = 0.02677843713276488

$ws.Range("I9").Value = 0.1062533062835484
$ws.Range("J9").Value = 0.1062533062835484
$ws.Range("M9").Value = 65.761079
$ws.Range("N9").Value = 197.283237
$ws.Range("O9").Value = 0.1293522275572212
$ws.Range("P9").Value = 0.1293522275572212
$ws.Range("Q9").Value = 3030.622166729773
$ws.Range("R9").Value = 27275.59950056795
$ws.Range("S9").Value = 0.01374410185309667
$ws.Range("T9").Value = 0.01374410185309668

$ws.Range("G10").Value = 41.187613
$ws.Range("H10").Value = 123.562839
$ws.Range("I10").Value = 0.09496120377532416
$ws.Range("J10").Value = 0.09496120377532417
$ws.Range("M10").Value = 145.7007446666667
$ws.Range("N10").Value = 437.1022340000001
$ws.Range("O10").Value = 0.2865937750105843
$ws.Range("P10").Value = 0.2865937750105843
$ws.Range("Q10").Value = 6001.065885142481
$ws.Range("R10").Value = 54009.59296628233
$ws.Range("S10").Value = 0.02721528986951949
$ws.Range("T10").Value = 0.0272152898695195

$ws.Range("G11").Value = 41.187613
$ws.Range("H11").Value = 123.562839
$ws.Range("I11").Value = 0.09496120377532416
$ws.Range("J11").Value = 0.09496120377532417
$ws.Range("O11").Value = 0.3320294904365841
$ws.Range("P11").Value = 0.3320294904365841
$ws.Range("Q11").Value = 6952.456827949735
$ws.Range("R11").Value = 62572.11145154762
$ws.Range("S11").Value = 0.03152992010076551
$ws.Range("T11").Value = 0.03152992010076552

$ws.Range("G12").Value = 41.187613
$ws.Range("H12").Value = 123.562839
$ws.Range("I12").Value = 0.09496120377532416
$ws.Range("J12").Value = 0.09496120377532417
$ws.Range("M12").Value = 128.1261546666667
$ws.Range("N12").Value = 384.378464
$ws.Range("O12").Value = 0.2520245069956105
$ws.Range("P12").Value = 0.2520245069956105
$ws.Range("Q12").Value = 5277.210473588811
$ws.Range("R12").Value = 47494.8942622993
$ws.Range("S12").Value = 0.02393255056518577
$ws.Range("T12").Value = 0.02393255056518578

$ws.Range("G13").Value = 41.187613
$ws.Range("H13").Value = 123.562839
$ws.Range("I13").Value = 0.09496120377532416
$ws.Range("J13").Value = 0.09496120377532417
$ws.Range("M13").Value = 65.761079
$ws.Range("N13").Value = 197.283237
$ws.Range("O13").Value = 0.1293522275572212
$ws.Range("P13").Value = 0.1293522275572212
$ws.Range("Q13").Value = 2708.541872314427
$ws.Range("R13").Value = 24376.87685082984
$ws.Range("S13").Value = 0.01228344323985339
$ws.Range("T13").Value = 0.01228344323985339

$ws.Range("G14").Value = 344.1819356666667
$ws.Range("H14").Value = 1032.545807
$ws.Range("I14").Value = 0.7935378757838636
$ws.Range("J14").Value = 0.7935378757838637
$ws.Range("M14").Value = 145.7007446666667
$ws.Range("N14").Value = 437.1022340000001
$ws.Range("O14").Value = 0.2865937750105843
$ws.Range("P14").Value = 0.2865937750105843
$ws.Range("Q14").Value = 50147.5643274481
$ws.Range("R14").Value = 451328.0789470329
$ws.Range("S14").Value = 0.2274230154347776
$ws.Range("T14").Value = 0.2274230154347776

$ws.Range("G15").Value = 344.1819356666667
$ws.Range("H15").Value = 1032.545807
$ws.Range("I15").Value = 0.7935378757838636
$ws.Range("J15").Value = 0.7935378757838637
$ws.Range("O15").Value = 0.3320294904365841
$ws.Range("P15").Value = 0.3320294904365841
$ws.Range("Q15").Value = 58097.80840377114
$ws.Range("R15").Value = 522880.2756339402
$ws.Range("S15").Value = 0.2634779765386456
$ws.Range("T15").Value = 0.2634779765386457

$ws.Range("G16").Value = 344.1819356666667
$ws.Range("H16").Value = 1032.545807
$ws.Range("I16").Value = 0.7935378757838636
$ws.Range("J16").Value = 0.7935378757838637
$ws.Range("M16").Value = 128.1261546666667
$ws.Range("N16").Value = 384.378464
$ws.Range("O16").Value = 0.2520245069956105
$ws.Range("P16").Value = 0.2520245069956105
$ws.Range("Q16").Value = 44098.70792270006
$ws.Range("R16").Value = 396888.3713043004
$ws.Range("S16").Value = 0.1999909919267722
$ws.Range("T16").Value = 0.1999909919267722

$ws.Range("G17").Value = 344.1819356666667
$ws.Range("H17").Value = 1032.545807
$ws.Range("I17").Value = 0.7935378757838636
$ws.Range("J17").Value = 0.7935378757838637
$ws.Range("M17").Value = 65.761079
$ws.Range("N17").Value = 197.283237
$ws.Range("O17").Value = 0.1293522275572212
$ws.Range("P17").Value = 0.1293522275572212
$ws.Range("Q17").Value = 22633.77546174858
$ws.Range("R17").Value = 203703.9791557372
$ws.Range("S17").Value = 0.1026458918836683
$ws.Range("T17").Value = 0.1026458918836683
